$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill column B (Attribute names) first, top to bottom ---
$ws.Range("B16").Value = "modX"
$ws.Range("B17").Value = "modY"
$ws.Range("B18").Value = "modXvalue"
$ws.Range("B19").Value = "modYvalue"
$ws.Range("B20").Value = "alignX"
$ws.Range("B21").Value = "alignY"

# --- Fill column C (Decorator / Required) next ---
$ws.Range("C16").Value = "IN => False"
$ws.Range("C17").Value = "IN => False"
$ws.Range("C18").Value = "IN => 100"
$ws.Range("C19").Value = "IN => 100"
$ws.Range("C20").Value = "IN => align.center"
$ws.Range("C21").Value = "IN => align.center"

# --- Fill column D (Type) next ---
$ws.Range("D16").Value = "bool"
$ws.Range("D17").Value = "bool"
$ws.Range("D18").Value = "int"
$ws.Range("D19").Value = "int"
$ws.Range("D20").Value = "String"

# --- Fill column E (Example Value) next ---
# E16/E17 reuse existing "true"/"false" text shared-strings (copy as values
# from E13/E14, which already hold them as text, to avoid Excel coercing the
# literal "true"/"false" into a boolean cell type)
$ws.Range("E13").Copy()
$ws.Range("E16").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E17").PasteSpecial(-4163)
$ws.Range("E18").Value = 90
$ws.Range("E19").Value = 90
$ws.Range("E20").Value = "align.center"
$ws.Range("E21").Value = "align.end"

# --- Fill column A (Widget name / section header) last ---
$ws.Range("A16").Value = "uPBOX"
$ws.Range("A6").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Rows.Item(16).RowHeight = 21

$ws.Application.CutCopyMode = $false

# --- Update dimension / selection ---
$ws.Range("A16").Select()
